$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.137.92'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '1.858.84'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'233.50"
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = "'0.4674"
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'0.2831"
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "'0.06459"
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = "'20.97"
$ws.Range('E10').Value = '  -3.59%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = "'0.07732"
$ws.Range('E11').Value = '  -3.61%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.858.85'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = "'93.59"
$ws.Range('E13').Value = '  -3.60%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'5.054"
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = "'0.6789"
$ws.Range('E15').Value = '  -0.96%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = "'266.22"
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '30.102.61'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = "'13.33"
$ws.Range('E18').Value = '  -4.57%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = "'0.000007560"
$ws.Range('E19').Value = '  -1.29%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = "'1.000"
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.103.62'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = "'5.153"
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').Value = "'6.095"
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').Value = "'9.280"
$ws.Range('E25').Value = '  -1.44%  '
$ws.Range('D26').Value = "'165.13"
$ws.Range('E26').Value = '  -2.04%  '
$ws.Range('D27').Value = "'18.48"
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('D28').Value = "'1.882"
$ws.Range('E28').Value = '  -3.60%  '
$ws.Range('D29').Value = "'1.372"
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = "'0.09843"
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('D32').Value = "'4.201"
$ws.Range('E32').Value = '  -3.97%  '
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('D34').Value = "'0.04650"
$ws.Range('D35').Value = "'1.113"
$ws.Range('E35').Value = '  -1.93%  '
$ws.Range('D36').Value = "'0.6849"
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('D38').Value = "'0.01811"
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('E39').Value = '  +3.58%  '
$ws.Range('D40').Value = "'6.281"
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').Value = "'70.57"
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('D42').Value = "'1.000"
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = "'0.8317"
$ws.Range('E43').Value = '  -1.16%  '
$ws.Range('D44').Value = "'1.877"
$ws.Range('E44').Value = '  -4.18%  '
$ws.Range('D45').Value = "'101.88"
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('D46').Value = "'0.4038"
$ws.Range('E46').Value = '  -3.08%  '
$ws.Range('D47').Value = "'9.138"
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').Value = "'923.71"
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = "'6.934"
$ws.Range('D50').Value = "'34.06"
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('D51').Value = "'0.05556"
$ws.Range('E51').Value = '  -2.25%  '
